$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add a new "2022-Q1" worksheet, positioned right before "总计", by
#    copying the "2021-Q4" sheet (same column layout / header wording
#    as the new sheet) and then overwriting its data rows.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($totalSheet)

# The freshly copied sheet is inserted immediately before "总计" and is
# named "2021-Q4 (2)" (or similar) by Excel - grab it via its position.
# (Re-fetch "总计" - any handle obtained before the copy has a stale
# .Index - then look one slot to its left.)
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# Make room for 4 data rows (the template only carried 1).
$newSheet.Rows.Item(3).Insert()
$newSheet.Rows.Item(3).Insert()
$newSheet.Rows.Item(3).Insert()

# Copy the formatting of the first data row down onto the 3 new rows.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A3:H5").PasteSpecial(-4122)

# Columns B (fund code) through G (ratios) must all stay TEXT, even
# though most of the values look numeric (leading zeros on fund codes,
# decimal-looking percentages) - matches every other quarter sheet.
$newSheet.Range("B2:G5").NumberFormat = "@"

# --- row 2 : 009562 ------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "009562"
$newSheet.Range("C2").Value = "工银瑞信中国机会全球配置股票(QDII)美元"
$newSheet.Range("D2").Value = "6.65"
$newSheet.Range("E2").Value = "92.85"
$newSheet.Range("F2").Value = "1.41"
$newSheet.Range("G2").Value = "0.0938"
$newSheet.Range("H2").Value = 9

# --- row 3 : 486001 ------------------------------------------------
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "486001"
$newSheet.Range("C3").Value = "工银瑞信中国机会全球配置股票(QDII)"
$newSheet.Range("D3").Value = "6.65"
$newSheet.Range("E3").Value = "92.85"
$newSheet.Range("F3").Value = "1.41"
$newSheet.Range("G3").Value = "0.0938"
$newSheet.Range("H3").Value = 9

# --- row 4 : 009563 ------------------------------------------------
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "009563"
$newSheet.Range("C4").Value = "工银瑞信中国机会全球配置股票(QDII)港币"
$newSheet.Range("D4").Value = "6.65"
$newSheet.Range("E4").Value = "92.85"
$newSheet.Range("F4").Value = "1.41"
$newSheet.Range("G4").Value = "0.0938"
$newSheet.Range("H4").Value = 9

# --- row 5 : 486002 ------------------------------------------------
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "486002"
$newSheet.Range("C5").Value = "工银全球精选股票(QDII)"
$newSheet.Range("D5").Value = "4.23"
$newSheet.Range("E5").Value = "94.60"
$newSheet.Range("F5").Value = "1.98"
$newSheet.Range("G5").Value = "0.0838"
$newSheet.Range("H5").Value = 5

# Now that every text value is safely stored, drop back to the default
# "Normal" style so B:G carry no explicit numFmt (matches every other
# quarter sheet - only A/H and the header row keep an explicit style).
$newSheet.Range("B2:G5").Style = "Normal"

# ---------------------------------------------------------------------
# 2. Update the "总计" overview sheet: insert a new top data row with
#    the 2022-Q1 summary, pushing the existing rows down.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Clean any stray formatting Insert() may have carried into B2:D2.
$totalSheet.Range("B2:D2").ClearFormats()

# Give A2 the same bold/bordered look used by every other index cell.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.37

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

Write-Host "2022-Q1 sheet added; summary sheet updated"
